$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# Locate the "Table 2" graphicFrame (shape id 3, name "Table 2") and update the
# "Domain Expiry" row / "Recommended Threshold" column cell: "< 3 months" -> "> 3 months".
for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $shape = $s.Shapes.Item($i)
    if ($shape.HasTable) {
        $table = $shape.Table
        for ($r = 1; $r -le $table.Rows.Count; $r++) {
            for ($c = 1; $c -le $table.Columns.Count; $c++) {
                $cellShape = $table.Cell($r, $c).Shape
                $tr = $cellShape.TextFrame.TextRange
                if ($tr.Text -eq "< 3 months") {
                    $tr.Text = "> 3 months"
                }
            }
        }
    }
}
